# Scheduled-runner price/profit refresh for Anima_Profits (Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW).
# Updates cached market-board columns H:N (currentAveragePrice*, LevePrice*, LeveProfit*) row by row.
# Some rows gain or lose an M (LeveProfitNQ) / N (LeveProfitHQ) cell entirely -- an empty-string
# assignment clears/removes the cell the same way the upstream export does.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1232.0588
$ws.Range("I32").Value = 1518
$ws.Range("J32").Value = 1112.9166
$ws.Range("K32").Value = 1518
$ws.Range("L32").Value = 1112.9166
$ws.Range("M32").Value = -1192
$ws.Range("N32").Value = -1764.9166
# Row 103
$ws.Range("H103").Value = 100607.8
$ws.Range("I103").Value = 111697.555
$ws.Range("K103").Value = 335092.665
$ws.Range("M103").Value = -334506.665
# Row 138
$ws.Range("H138").Value = 192697.17
$ws.Range("J138").Value = 295681.66
$ws.Range("L138").Value = 887044.98
$ws.Range("N138").Value = -897324.98
# Row 141
$ws.Range("H141").Value = 5775
$ws.Range("I141").Value = 2798.889
$ws.Range("J141").Value = 12471.25
$ws.Range("K141").Value = 8396.667000000001
$ws.Range("L141").Value = 37413.75
$ws.Range("M141").Value = -3216.667000000001
$ws.Range("N141").Value = -47773.75

$ws = $wb.Worksheets.Item("ARM")
# Row 52
$ws.Range("H52").Value = 82520
$ws.Range("J52").Value = 82520
$ws.Range("L52").Value = 82520
$ws.Range("N52").Value = -83156
# Row 61
$ws.Range("H61").Value = 2964.5
$ws.Range("I61").Value = 2705.625
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2705.625
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2493.625
$ws.Range("N61").Value = -4424
# Row 63
$ws.Range("H63").Value = 6014.077
$ws.Range("I63").Value = 3259.4
$ws.Range("K63").Value = 3259.4
$ws.Range("M63").Value = -2573.4
# Row 66
$ws.Range("H66").Value = 6014.077
$ws.Range("I66").Value = 3259.4
$ws.Range("K66").Value = 16297
$ws.Range("M66").Value = -12865
# Row 102
$ws.Range("H102").Value = 1920
$ws.Range("I102").Value = 1975
$ws.Range("J102").Value = 1700
$ws.Range("K102").Value = 1975
$ws.Range("L102").Value = 1700
$ws.Range("M102").Value = -353
$ws.Range("N102").Value = -4944
# Row 132
$ws.Range("H132").Value = 4743.2
$ws.Range("I132").Value = 4066.6296
$ws.Range("J132").Value = 6148.385
$ws.Range("K132").Value = 12199.8888
$ws.Range("L132").Value = 18445.155
$ws.Range("M132").Value = -9669.888800000001
$ws.Range("N132").Value = -23505.155
# Row 136
$ws.Range("H136").Value = 2964.5
$ws.Range("I136").Value = 2705.625
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 8116.875
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -5566.875
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 62501788
$ws.Range("I86").Value = 76924820
$ws.Range("J86").Value = 1996.6666
$ws.Range("K86").Value = 76924820
$ws.Range("L86").Value = 1996.6666
$ws.Range("M86").Value = -76923697
$ws.Range("N86").Value = -4242.6666
# Row 89
$ws.Range("H89").Value = 62501788
$ws.Range("I89").Value = 76924820
$ws.Range("J89").Value = 1996.6666
$ws.Range("K89").Value = 384624100
$ws.Range("L89").Value = 9983.333000000001
$ws.Range("M89").Value = -384618484
$ws.Range("N89").Value = -21215.333

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1437.7916
$ws.Range("I31").Value = 871.4857
$ws.Range("J31").Value = 2962.4614
$ws.Range("K31").Value = 871.4857
$ws.Range("L31").Value = 2962.4614
$ws.Range("M31").Value = -576.4857
$ws.Range("N31").Value = -3552.4614
# Row 34
$ws.Range("H34").Value = 1437.7916
$ws.Range("I34").Value = 871.4857
$ws.Range("J34").Value = 2962.4614
$ws.Range("K34").Value = 871.4857
$ws.Range("L34").Value = 2962.4614
$ws.Range("M34").Value = -669.4857
$ws.Range("N34").Value = -3366.4614
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 13890889
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 13890889
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 41672667
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -41674289
# Row 72
$ws.Range("H72").Value = 13890889
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 13890889
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 125018001
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -125026113
# Row 74
$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = ""
# Row 77
$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = ""
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = ""
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1678
$ws.Range("I102").Value = 1570.6666
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1570.6666
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 51.33339999999998
$ws.Range("N102").Value = -5244
# Row 132
$ws.Range("H132").Value = 2230.652
$ws.Range("I132").Value = 1579.5333
$ws.Range("K132").Value = 4738.5999
$ws.Range("M132").Value = -2208.5999

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4503.6333
$ws.Range("I22").Value = 790.7143
$ws.Range("J22").Value = 7752.4375
$ws.Range("K22").Value = 790.7143
$ws.Range("L22").Value = 7752.4375
$ws.Range("M22").Value = -495.7143
$ws.Range("N22").Value = -8342.4375
# Row 27
$ws.Range("H27").Value = 4503.6333
$ws.Range("I27").Value = 790.7143
$ws.Range("J27").Value = 7752.4375
$ws.Range("K27").Value = 790.7143
$ws.Range("L27").Value = 7752.4375
$ws.Range("M27").Value = -683.7143
$ws.Range("N27").Value = -7966.4375
